$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.673.49"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "1.637.80"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "213.16"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").Value = "19.21"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("D11").Value = "0.0841"
$ws.Range("E11").Value = "  +3.18%  "
$ws.Range("D13").Value = "1.629.88"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "26.679.11"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "63.39"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").Value = "218.81"
$ws.Range("E19").Value = "  +7.67%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "4.31"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "9.50"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "148.59"
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "6.86"
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  +4.00%  "
$ws.Range("D33").Value = "2.95"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").Value = "1.196.14"
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("E37").Value = "  +6.03%  "
$ws.Range("D38").Value = "0.810"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.795"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.41"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("D44").Value = "1.774.49"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").Value = "54.83"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D49").Value = "7.66"
$ws.Range("E49").Value = "  +5.46%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  +0.19%  "
